$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 162192
$ws.Range("C4").Value = 153201
$ws.Range("C7").Value = 5.54
$ws.Range("C8").Value = 64.72
